$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 63057
$ws.Range("I33").Value = 83690.336
$ws.Range("J33").Value = 1157
$ws.Range("K33").Value = 83690.336
$ws.Range("L33").Value = 1157
$ws.Range("M33").Value = -83461.336
$ws.Range("N33").Value = -1615

$ws.Range("H69").Value = 10749.667
$ws.Range("I69").Value = 9999.666999999999
$ws.Range("K69").Value = 29999.001
$ws.Range("M69").Value = -29125.001

$ws.Range("H70").Value = 852619.0600000001
$ws.Range("I70").Value = 1701865
$ws.Range("J70").Value = 3373.1667
$ws.Range("K70").Value = 5105595
$ws.Range("L70").Value = 10119.5001
$ws.Range("M70").Value = -5105325
$ws.Range("N70").Value = -10659.5001

$ws.Range("H72").Value = 10749.667
$ws.Range("I72").Value = 9999.666999999999
$ws.Range("K72").Value = 89997.003
$ws.Range("M72").Value = -85629.003

$ws.Range("H73").Value = 852619.0600000001
$ws.Range("I73").Value = 1701865
$ws.Range("J73").Value = 3373.1667
$ws.Range("K73").Value = 5105595
$ws.Range("L73").Value = 10119.5001
$ws.Range("M73").Value = -5104659
$ws.Range("N73").Value = -11991.5001

$ws.Range("H82").Value = 953.625
$ws.Range("I82").Value = 375.57144
$ws.Range("K82").Value = 1126.71432
$ws.Range("M82").Value = -720.71432

$ws.Range("H85").Value = 953.625
$ws.Range("I85").Value = 375.57144
$ws.Range("K85").Value = 1126.71432
$ws.Range("M85").Value = 277.28568

$ws.Range("H94").Value = 2087.5
$ws.Range("I94").Value = 2087.5
$ws.Range("K94").Value = 2087.5
$ws.Range("M94").Value = -1636.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2989.131
$ws.Range("I32").Value = 1595.7885
$ws.Range("K32").Value = 1595.7885
$ws.Range("M32").Value = -1308.7885

$ws.Range("H45").Value = 1923.5
$ws.Range("I45").Value = 1734.7273
$ws.Range("J45").Value = 4000
$ws.Range("K45").Value = 1734.7273
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = -1357.7273
$ws.Range("N45").Value = -4754

$ws.Range("H61").Value = 55558396
$ws.Range("I61").Value = 83334216
$ws.Range("K61").Value = 83334216
$ws.Range("M61").Value = -83334004

$ws.Range("H64").Value = 49995
$ws.Range("J64").Value = 49995
$ws.Range("L64").Value = 49995
$ws.Range("N64").Value = -50491

$ws.Range("H67").Value = 49995
$ws.Range("J67").Value = 49995
$ws.Range("L67").Value = 49995
$ws.Range("N67").Value = -51711

$ws.Range("H88").Value = 15153797
$ws.Range("I88").Value = 55556856
$ws.Range("J88").Value = 2649.875
$ws.Range("K88").Value = 55556856
$ws.Range("L88").Value = 2649.875
$ws.Range("M88").Value = -55556450
$ws.Range("N88").Value = -3461.875

$ws.Range("H91").Value = 15153797
$ws.Range("I91").Value = 55556856
$ws.Range("J91").Value = 2649.875
$ws.Range("K91").Value = 55556856
$ws.Range("L91").Value = 2649.875
$ws.Range("M91").Value = -55555452
$ws.Range("N91").Value = -5457.875

$ws.Range("H132").Value = 27779432
$ws.Range("I132").Value = 31251502
$ws.Range("K132").Value = 93754506
$ws.Range("M132").Value = -93751976

$ws.Range("H136").Value = 55558396
$ws.Range("I136").Value = 83334216
$ws.Range("K136").Value = 250002648
$ws.Range("M136").Value = -250000098

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 12749.5
$ws.Range("I49").Value = 5499
$ws.Range("J49").Value = 20000
$ws.Range("K49").Value = 5499
$ws.Range("L49").Value = 20000
$ws.Range("M49").Value = -5260
$ws.Range("N49").Value = -20478

$ws.Range("H105").Value = 1885.875
$ws.Range("I105").Value = 1774.9231
$ws.Range("K105").Value = 1774.9231
$ws.Range("M105").Value = -27.92309999999998

$ws.Range("H134").Value = 3260.5757
$ws.Range("I134").Value = 3212.9033
$ws.Range("K134").Value = 9638.7099
$ws.Range("M134").Value = -7103.7099

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2991.6377
$ws.Range("I31").Value = 1456.8889
$ws.Range("K31").Value = 1456.8889
$ws.Range("M31").Value = -1161.8889

$ws.Range("H34").Value = 2991.6377
$ws.Range("I34").Value = 1456.8889
$ws.Range("K34").Value = 1456.8889
$ws.Range("M34").Value = -1254.8889

$ws.Range("H58").Value = 2039.6111
$ws.Range("I58").Value = 1907.3846
$ws.Range("J58").Value = 2383.4
$ws.Range("K58").Value = 1907.3846
$ws.Range("L58").Value = 2383.4
$ws.Range("M58").Value = -1704.3846
$ws.Range("N58").Value = -2789.4

$ws.Range("H132").Value = 2829.8572
$ws.Range("I132").Value = 2807.6667
$ws.Range("K132").Value = 8423.000100000001
$ws.Range("M132").Value = -5893.000100000001

$ws.Range("H136").Value = 2039.6111
$ws.Range("I136").Value = 1907.3846
$ws.Range("J136").Value = 2383.4
$ws.Range("K136").Value = 5722.1538
$ws.Range("L136").Value = 7150.200000000001
$ws.Range("M136").Value = -3172.1538
$ws.Range("N136").Value = -12250.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 9513183
$ws.Range("I4").Value = 10233565
$ws.Range("K4").Value = 30700695
$ws.Range("M4").Value = -30700583

$ws.Range("H34").Value = 1642.0476
$ws.Range("I34").Value = 175
$ws.Range("J34").Value = 1796.4736
$ws.Range("K34").Value = 525
$ws.Range("L34").Value = 5389.4208
$ws.Range("M34").Value = -441
$ws.Range("N34").Value = -5557.4208

$ws.Range("H40").Value = 269.8
$ws.Range("I40").Value = 49
$ws.Range("K40").Value = 196
$ws.Range("M40").Value = -127

$ws.Range("H75").Value = 412.69232
$ws.Range("I75").Value = 245
$ws.Range("K75").Value = 735
$ws.Range("M75").Value = 263

$ws.Range("H78").Value = 412.69232
$ws.Range("I78").Value = 245
$ws.Range("K78").Value = 2205
$ws.Range("M78").Value = 2787

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3143.353
$ws.Range("I80").Value = 3241.6667
$ws.Range("K80").Value = 3241.6667
$ws.Range("M80").Value = -2243.6667

$ws.Range("H83").Value = 3143.353
$ws.Range("I83").Value = 3241.6667
$ws.Range("K83").Value = 16208.3335
$ws.Range("M83").Value = -11216.3335

$ws.Range("H132").Value = 4467.225
$ws.Range("I132").Value = 4041.8572
$ws.Range("K132").Value = 12125.5716
$ws.Range("M132").Value = -9595.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1365.3334
$ws.Range("I82").Value = 1023.6667
$ws.Range("K82").Value = 1023.6667
$ws.Range("M82").Value = -662.6667

$ws.Range("H85").Value = 1365.3334
$ws.Range("I85").Value = 1023.6667
$ws.Range("K85").Value = 1023.6667
$ws.Range("M85").Value = 224.3333

$ws.Range("H132").Value = 7069.9443
$ws.Range("I132").Value = 3173.625
$ws.Range("J132").Value = 10187
$ws.Range("K132").Value = 9520.875
$ws.Range("L132").Value = 30561
$ws.Range("M132").Value = -6990.875
$ws.Range("N132").Value = -35621

$ws.Range("H136").Value = 3292.5
$ws.Range("I136").Value = 3199.5908
$ws.Range("K136").Value = 9598.7724
$ws.Range("M136").Value = -7048.7724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6671856
$ws.Range("I81").Value = 2314.6316
$ws.Range("K81").Value = 4629.2632
$ws.Range("M81").Value = -3568.2632

$ws.Range("H84").Value = 6671856
$ws.Range("I84").Value = 2314.6316
$ws.Range("K84").Value = 23146.316
$ws.Range("M84").Value = -17842.316

$ws.Range("H100").Value = 5231.8184
$ws.Range("I100").Value = 6906.25
$ws.Range("K100").Value = 13812.5
$ws.Range("M100").Value = -13271.5

$ws.Range("H132").Value = 4098.2256
$ws.Range("I132").Value = 4258.55
$ws.Range("K132").Value = 12775.65
$ws.Range("M132").Value = -10245.65

$ws.Range("H136").Value = 3899.8696
$ws.Range("J136").Value = 5590.077
$ws.Range("L136").Value = 16770.231
$ws.Range("N136").Value = -21870.231
